# edit.ps1
# Applies the "dark energy essay" -> "art appreciation essay" rewrite described
# by the commit diff, using Word COM-interop (Find/Replace + paragraph insert).

$d = $word.ActiveDocument

function Replace-InPara {
    param(
        [int]$ParaIndex,
        [string]$OldText,
        [string]$NewText
    )
    $p = $d.Paragraphs($ParaIndex)
    $ok = $p.Range.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if (-not $ok) {
        $preview = $OldText
        if ($preview.Length -gt 60) { $preview = $preview.Substring(0, 60) }
        Write-Output "WARNING: replacement not found (paragraph $ParaIndex): $preview"
    }
    return $ok
}

# ---------------------------------------------------------------------------
# Title
# ---------------------------------------------------------------------------
Replace-InPara 1 `
    "Unraveling the Enigma of Dark Energy: An Astronomical Odyssey" `
    "The Crucible of Creativity: Exploring the Essence of Art Appreciation"

# ---------------------------------------------------------------------------
# Author name line ("Dr. Sophia Adamson" -> "Ms. Isabella Rodriguez")
# ---------------------------------------------------------------------------
Replace-InPara 2 "Dr" "Ms"
Replace-InPara 2 " Sophia Adamson" " Isabella Rodriguez"

# ---------------------------------------------------------------------------
# Email line (sophia.adamson@stellarobservatory.org -> isabelarodriguezteacher@gmail.com)
# ---------------------------------------------------------------------------
Replace-InPara 3 "sophia" "isabelarodriguezteacher@gmail"
Replace-InPara 3 "adamson@stellarobservatory.org" "com"

# ---------------------------------------------------------------------------
# Body paragraph (paragraph 5) - three sentence groups separated by <br><br>
# ---------------------------------------------------------------------------

# --- group 1 ---
Replace-InPara 5 `
    "In the vast expanse of the cosmos, hidden behind the tapestry of stars and galaxies, lies a profound enigma that has captivated the minds of scientists and astronomers alike: dark energy" `
    "In the realm of human experience, art stands as a testament to our ingenuity and expressiveness"

Replace-InPara 5 `
    " This enigmatic force, believed to permeate the fabric of space, stands as one of the greatest mysteries in modern physics" `
    " It transcends linguistic boundaries, captivating souls with its visual symphony and emotional resonance"

Replace-InPara 5 `
    " It is an invisible entity that dominates the universe, driving its expansion at an ever-accelerating pace" `
    " Throughout history, renowned artists have employed various mediums, from paint and clay to music and literature, to weave narratives that mirror the tapestry of our shared existence"

Replace-InPara 5 `
    " Its existence challenges our understanding of gravity and the fundamental forces that govern the universe, beckoning us to embark on an astronomical odyssey to unravel its secrets" `
    " As educators, it is our privilege to unveil the secrets of art appreciation, illuminating the enigmatic connection between artwork and observer. Our journey through the vast landscape of artistry promises to unlock new perspectives and inspire boundless creativity within our students"

# --- group 2 ---
Replace-InPara 5 `
    "As we peer into the depths of the cosmos, distant supernovae reveal a peculiar tale: the expansion of the universe is not decelerating as expected under the influence of gravity, but rather accelerating" `
    "We shall embark on a quest to unravel the mysteries embedded within iconic paintings, tracing the evolution of artistic styles and techniques across different eras"

Replace-InPara 5 `
    " This counterintuitive observation, first hinted at by Edwin Hubble in the 1920s and later confirmed by subsequent studies, has profound implications" `
    " Each masterpiece invites us to decipher its symbolism, delve into its historical context, and appreciate its unique contribution to the global tapestry of human creativity"

Replace-InPara 5 `
    " It suggests the existence of a mysterious energy permeating space that counteracts the pull of gravity and propels the universe's expansion. This enigmatic entity is what we refer to as dark energy" `
    " From the sublime landscapes of the Renaissance to the bold strokes of modernism, we will trace the threads that bind artists and their audiences, exploring how art reflects and shapes societal norms, values, and beliefs"

# --- group 3 ---
Replace-InPara 5 `
    "The nature of dark energy remains elusive, shrouded in uncertainty" `
    "Further, we will delve into the enchanting world of music, where melodies and harmonies dance in perfect unison"

Replace-InPara 5 `
    " Scientists have proposed various theories to explain its existence, ranging from modifications to Einstein's theory of gravity to the presence of a cosmological constant, a constant energy density permeating the universe" `
    " Music has the power to transcend languages and cultures, stirring emotions that words cannot express"

Replace-InPara 5 `
    " Other hypotheses include scalar fields, dynamical dark energy models, and even the possibility of a multiverse, where dark energy arises from the interactions between parallel universes" `
    " We will explore the intricacies of musical composition, dissecting the interplay of rhythm, pitch, and timbre, and trace the evolution of genres from classical to contemporary"

Replace-InPara 5 `
    " Yet, despite these theoretical advances, the true identity of dark energy remains hidden, waiting to be unveiled by future observations and experiments" `
    " By attuning our ears to the symphony of sound, we will cultivate a deeper understanding of the emotional and intellectual power of music"

# ---------------------------------------------------------------------------
# Summary paragraph (paragraph 7)
# ---------------------------------------------------------------------------
Replace-InPara 7 `
    "The quest to understand dark energy is a captivating journey at the frontiers of physics and astronomy" `
    "Through exploration of diverse artistic mediums, this essay has illuminated the essence of art appreciation, revealing the profound interconnectedness between artwork and observer"

Replace-InPara 7 `
    " This mysterious force, responsible for the accelerated expansion of the universe, challenges our current understanding of gravity and the fundamental forces that govern the " `
    " By "

Replace-InPara 7 `
    "cosmos" `
    "unveiling the mysteries of art history, tracing the evolution of styles, and exploring the emotional and intellectual power of music, it has provided a framework for fostering a lifelong appreciation for artistic endeavors"

Replace-InPara 7 `
    " As we delve deeper into the mysteries of dark energy, we embark on an odyssey of discovery, seeking to unravel the enigmas that lie hidden within the vast expanses of space. Through observations, experiments, and theoretical explorations, we strive to shed light on this enigmatic force, unlocking the secrets that hold the key to comprehending the ultimate fate of our universe" `
    " Encouraging students to engage with art critically and creatively empowers them to become informed and thoughtful consumers of culture, preparing them to navigate the complex tapestry of human expression that surrounds us"

# ---------------------------------------------------------------------------
# New trailing empty paragraph after the Summary paragraph
# ---------------------------------------------------------------------------
$d.Paragraphs.Add() | Out-Null

Write-Output "Done."
